$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45070
$ws.Range("J2").Value = 800
$ws.Range("D3").Value = 45020
$ws.Range("J3").Value = 1200
$ws.Range("D4").Value = 44883
$ws.Range("J4").Value = 500
$ws.Range("D5").Value = 44910
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("P5").Value = 633
$ws.Range("D6").Value = 44848
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1750
$ws.Range("P6").Value = 583
$ws.Range("D7").Value = 44953
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2250
$ws.Range("P7").Value = 750
$ws.Range("D8").Value = 45007
$ws.Range("J8").Value = 1160
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2250
$ws.Range("P8").Value = 750
$ws.Range("D9").Value = 44827
$ws.Range("J9").Value = 1200
$ws.Range("D10").Value = 44992
$ws.Range("J10").Value = 1040
$ws.Range("D11").Value = 45013
$ws.Range("J11").Value = 1100
$ws.Range("D12").Value = 44971
$ws.Range("J12").Value = 1000
$ws.Range("D13").Value = 45084
$ws.Range("J13").Value = 900
$ws.Range("D14").Value = 44970
$ws.Range("J14").Value = 800
$ws.Range("D15").Value = 45041
$ws.Range("J15").Value = 1160
$ws.Range("D16").Value = 45034
$ws.Range("J16").Value = 1100
$ws.Range("D17").Value = 45006
$ws.Range("D18").Value = 45077
$ws.Range("J18").Value = 760
$ws.Range("D19").Value = 44999
$ws.Range("D20").Value = 44951
$ws.Range("J20").Value = 800
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2250
$ws.Range("P20").Value = 750
$ws.Range("D21").Value = 44964
$ws.Range("J21").Value = 1000
$ws.Range("D22").Value = 44881
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1900
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 1950
$ws.Range("P22").Value = 650
$ws.Range("D23").Value = 44985
$ws.Range("D24").Value = 44685
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 1500
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = 1750
$ws.Range("P24").Value = 583
$ws.Range("D25").Value = 45091
$ws.Range("J25").Value = 800
$ws.Range("D26").Value = 45062
$ws.Range("J26").Value = 1100
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 2500
$ws.Range("M26").Value = 2250
$ws.Range("P26").Value = 750
$ws.Range("D27").Value = 44965
$ws.Range("J27").Value = 1120
$ws.Range("D28").Value = 44911
$ws.Range("J28").Value = 700
$ws.Range("K28").Value = 1800
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = 1900
$ws.Range("P28").Value = 633
$ws.Range("D29").Value = 45028
$ws.Range("J29").Value = 1000
$ws.Range("D30").Value = 45035
$ws.Range("J30").Value = 1100
$ws.Range("K30").Value = 2000
$ws.Range("L30").Value = 2500
$ws.Range("M30").Value = 2250
$ws.Range("P30").Value = 750
$ws.Range("D31").Value = 44978
$ws.Range("J31").Value = 1000
